# JLV:25.08.2022: Descarga de documentos Word
#
# Adds the new batch of document IDs pasted into Sheet1 column A (rows
# 3-51, right under the existing header/first ID), which pushes Excel's
# "duplicate values" conditional formatting to re-span the column, and
# finally leaves the active cell on C12 (where the user clicked next).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ----------------------------------------------------------------------
# 1) New IDs, rows 3-19 -- these were entered/pasted as TEXT (shared
#    strings) in the source workbook.
# ----------------------------------------------------------------------
$textValues = [ordered]@{
    3  = "88062"
    4  = "12561"
    5  = "10444"
    6  = "12434"
    7  = "12711"
    8  = "12801"
    9  = "15458"
    10 = "15756"
    11 = "47001"
    12 = "47816"
    13 = "52771"
    14 = "17444"
    15 = "35503"
    16 = "39744"
    17 = "77891"
    18 = "57004"
    19 = "57007"
}
foreach ($r in $textValues.Keys) {
    $ws.Cells.Item($r, 1).Value2 = $textValues[$r]
}
# Pick up the same look (style) as the rest of the ID column.
$ws.Range("A2").Copy()
$ws.Range("A3:A19").PasteSpecial(-4122)

# ----------------------------------------------------------------------
# 2) New IDs, rows 20-51 -- entered as plain NUMBERS. Apply the column's
#    look *before* writing the values (instead of after, as above) so the
#    destination cells already carry a "General" numeric format when the
#    value lands -- that keeps them real numbers instead of text, and
#    without it Excel would otherwise mint a redundant duplicate style.
# ----------------------------------------------------------------------
$numValues = [ordered]@{
    20 = 90964
    21 = 7750
    22 = 77657
    23 = 8281
    24 = 15864
    25 = 60426
    26 = 62164
    27 = 22990
    28 = 11062
    29 = 37791
    30 = 8959
    31 = 23361
    32 = 22521
    33 = 11839
    34 = 39740
    35 = 8161
    36 = 89616
    37 = 54735
    38 = 67231
    39 = 75877
    40 = 80064
    41 = 82856
    42 = 86180
    43 = 88743
    44 = 89090
    45 = 89344
    46 = 93493
    47 = 66968
    48 = 48433
    49 = 62811
    50 = 69223
    51 = 90965
}
$ws.Range("A2").Copy()
$ws.Range("A20:A51").PasteSpecial(-4122)
foreach ($r in $numValues.Keys) {
    $ws.Cells.Item($r, 1).Value = $numValues[$r]
}

# ----------------------------------------------------------------------
# 3) The "duplicate values" conditional formatting on column A used to be
#    split into several adjacent blocks (A2:A8 / A9:A26 / A27:A30 /
#    A31:A177). Typing the new rows in the middle of that range makes
#    Excel fold the now-contiguous blocks back together, leaving just
#    A2:A51 and A52:A177.
# ----------------------------------------------------------------------
$allConditions = $ws.Cells.FormatConditions
for ($i = $allConditions.Count; $i -ge 1; $i--) {
    $fc = $allConditions.Item($i)
    $addr = $fc.AppliesTo.Address()
    if ($addr -eq "`$A`$9:`$A`$26") {
        $fc.Delete()
    }
    elseif ($addr -eq "`$A`$27:`$A`$30") {
        $fc.Delete()
    }
    elseif ($addr -eq "`$A`$31:`$A`$177") {
        $fc.ModifyAppliesToRange($ws.Range("A52:A177"))
    }
    elseif ($addr -eq "`$A`$2:`$A`$8") {
        $fc.ModifyAppliesToRange($ws.Range("A2:A51"))
    }
}

# ----------------------------------------------------------------------
# 4) Leave the selection where the user clicked next.
# ----------------------------------------------------------------------
$ws.Range("C12").Select()
